# [FIX] added missing/recently added boardGame fields to export template
#
# Inserts a "Popis" (${description}) column right after "Nazov Hry", and
# appends four new trailing columns - "Vek" (${ageRange}),
# "Kooperativna" (${isCooperative}), "Rozsirenie" (${isExtension}) and
# "Jednorazova" (${isOneTimePlay}) - before the existing "Autori" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert the new "Popis" column at C ------------------------
# This shifts old C:F (Min./Max. players, play time, Autori) one column
# right, to D:G, carrying their widths/styles/values along with them.
$ws.Range("C1:C2").EntireColumn.Insert()

# --- Step 2: insert 4 blank columns before the (now) "Autori" column ---
# Old "Autori" is now at G; inserting G:J pushes it to K, again carrying
# its width/style/value with it.
$ws.Range("G1:J2").EntireColumn.Insert()

# --- Step 3: fill in the new header row (row 1) -------------------------
$ws.Range("C1").Value = "Popis"
$ws.Range("G1").Value = "Vek"
$ws.Range("H1").Value = "Kooperatívna"
$ws.Range("I1").Value = "Rozšírenie"
$ws.Range("J1").Value = "Jednorázová"

# --- Step 4: fill in the new placeholder row (row 2) ---------------------
$ws.Range("C2").Value = '${description}'
$ws.Range("G2").Value = '${ageRange}'
$ws.Range("H2").Value = '${isCooperative}'
$ws.Range("I2").Value = '${isExtension}'
$ws.Range("J2").Value = '${isOneTimePlay}'

# --- Step 5: size the newly inserted columns -----------------------------
# (A:B and the shifted D:G/K keep their original widths automatically.)
$ws.Columns.Item(3).ColumnWidth = 50.1    # C  - Popis            (~50.99)
$ws.Columns.Item(7).ColumnWidth = 14.6    # G  - Vek              (~15.56)
$ws.Columns.Item(8).ColumnWidth = 12      # H  - Kooperativna     (~12.91)
$ws.Columns.Item(9).ColumnWidth = 12      # I  - Rozsirenie       (~12.91)
$ws.Columns.Item(10).ColumnWidth = 12     # J  - Jednorazova      (~12.91)

Write-Output "BoardGamesTemplate updated: added Popis/Vek/Kooperativna/Rozsirenie/Jednorazova columns"
